$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The label in A11 was "Sunny Days"; update it to "Holidays" (removing the
# now-unused "Sunny Days" shared string).
$ws.Range("A11").Value = "Holidays"

# Update the selected cell on the sheet to A12.
$ws.Range("A12").Select()
